$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The figures in column C (cfu counts) were recalculated - every value is
# multiplied by 4 (dilution factor correction).
for ($r = 2; $r -le 33; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value()
    $cell.Value = $current * 4
}

# Row 2 previously held an oversized helper cell (F2) with a big 32pt font
# and a taller row (ht=40). That stray cell/formatting is removed so the
# row goes back to the sheet's normal height and the used range shrinks
# back down to column D.
$ws.Range("F2").Clear()
$ws.Rows.Item(2).AutoFit()

# Column F is given an explicit width (as if Excel had just best-fit it)
# even though the sheet's data no longer extends into it.
$ws.Columns.Item(6).ColumnWidth = 22.6667

# Update the active selection to F3, matching where the user clicked next.
$ws.Range("F3").Select()

# Reposition the workbook window (best effort - matches author's window move).
$excel.ActiveWindow.Left = 10520
